# "updated import spreadhseet with errors for presentation"
# Simulate missing/erroneous source data by blanking out a handful of cells
# in the gift-drive import sheet: first_name (A8, A16, A31), gender (C8, C14),
# size_pants (D8, D13) and size_shirt (E31). Excel will drop the now-unused
# shared strings ("Philip", "Isabella", "Jasmin") from the shared-string
# table automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @("A8", "C8", "D8", "D13", "C14", "A16", "A31", "E31")
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
